# Update the metadata sheet with the finalized "Raw Data Citation" and
# "Variables Generated in Dataset:" text (the placeholder citation date /
# anonymized author name are replaced with the real retrieval date and
# author, and the variable list is expanded with the new Crab/Fish/
# NonShrimpInvert CPUE & MPUE metrics), then move the active selection to
# B10 to match the author's final cursor position when the file was saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B8").Value = 'SEAMAP-SA Data Management Work Group . 2022, November, 1. SEAMAP-SA online database. Retrieved from: http://www.dnr.sc.gov/SEAMAP/data.html by: Lela Schlenker.'
$ws.Range("B10").Value = 'Summer and Fall: CPUE (catch per unit effort), MPUE (biomass per unit effort), CrabCPUE (catch per unit effort of macroinvertebrates), CrabMPUE (biomass per unit effort of macroinvertebrates), FishCPUE (catch per unit effort of finfish), FishMPUE (biomass per unit effort of finfish),  Species Richness, DtoPBioMRatio (demersal to pelagic biomass ratio), shannon wiener diversity index, NonShrimpInvertCPUE (catch per unit effort of non-shrimp macroinvertebrates)'

$ws.Range("B10").Select()
